# Apartado 3: criterios de comparación. FINALIZADO
#
# Update the "responsable" column (D) for a few rows of the TG2 repartition
# sheet: José joins the "Criterios de comparación" / "Evaluación criterios
# tecnología 1" tasks, and is now also assigned to "Situación 1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 -> "4.1) Evaluación criterios tecnología 1"
$ws.Range("D9").Value = "José y Fernando"

# Row 13 -> "6.1) Situación 1" (previously empty)
$ws.Range("D13").Value = "José"

# Row 7 -> "3) Criterios de comparación"
$ws.Range("D7").Value = "José, Fernando Criterios E y F"

# Move the active selection to D5, matching the saved view state.
[void]$ws.Range("D5").Select()
